# Update the "Input" sheet with the new PRO number, instrument SN, and
# clear the ICA SN value, then force a recalculation so the dependent
# formulas on "Template_printout" (and the Cell# computation on Input)
# pick up the new results.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")

$ws.Range("B3").Value = 20581576
$ws.Range("B4").Value = "A01606"
$ws.Range("B5").Value = ""

$excel.Calculate()
